$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F334" = 195995
    "F335" = 130923
    "G335" = 3001
    "F336" = 101862
    "G336" = 3378
    "F337" = 103905
    "G337" = 2959
    "F338" = 226837
    "G338" = 3186
    "F339" = 657739
    "G339" = 5487
    "F340" = 382342
    "G340" = 3294
    "F342" = 179404
    "G342" = 3068
    "F343" = 132430
    "G343" = 2967
    "F344" = 135448
    "G344" = 2487
    "F345" = 291612
    "F346" = 670176
    "G346" = 4785
    "F347" = 342310
    "G347" = 2912
    "F348" = 232031
    "G348" = 3251
    "F349" = 159061
    "G349" = 2752
    "F350" = 127010
    "G350" = 2788
    "F351" = 150439
    "G351" = 2823
    "F352" = 306784
    "G352" = 3538
    "F353" = 721442
    "G353" = 5255
    "F355" = 221970
    "G355" = 3450
    "F356" = 159903
    "G356" = 2873
    "F357" = 138435
    "G357" = 3025
    "F358" = 157108
    "G358" = 2597
    "F359" = 319970
    "G359" = 3334
    "F360" = 745871
    "G360" = 5108
    "F363" = 187876
    "G363" = 2756
    "F365" = 179992
    "G365" = 2349
    "F367" = 761386
    "F369" = 233995
    "G369" = 2572
    "F370" = 181892
    "F373" = 343449
    "G373" = 2343
    "F374" = 762705
    "G374" = 3361
    "F375" = 349893
    "F376" = 218839
    "G376" = 2187
    "F377" = 174488
    "F378" = 153864
    "G378" = 1507
    "F379" = 173809
    "G379" = 1575
    "F380" = 332973
    "G380" = 1946
    "F381" = 702566
    "G381" = 3088
    "F382" = 335228
    "G382" = 2053
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
